$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.878.39'
$ws.Range('E2').Value = '  -3.46%  '
$ws.Range('D3').Value = '2.532.57'
$ws.Range('E3').Value = '  -4.73%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'511.07"
$ws.Range('E5').Value = '  -2.47%  '
$ws.Range('D6').Value = "'139.06"
$ws.Range('E6').Value = '  -3.79%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -2.98%  '
$ws.Range('D9').Value = "'6.47"
$ws.Range('E9').Value = '  -7.97%  '
$ws.Range('D10').Value = "'0.0990"
$ws.Range('E10').Value = '  -3.50%  '
$ws.Range('E11').Value = '  -3.57%  '
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('D13').Value = '2.979.13'
$ws.Range('E13').Value = '  -4.55%  '
$ws.Range('D14').Value = '56.914.19'
$ws.Range('E14').Value = '  -3.43%  '
$ws.Range('D15').Value = "'19.96"
$ws.Range('E15').Value = '  -5.33%  '
$ws.Range('E16').Value = '  -3.20%  '
$ws.Range('D17').Value = '2.494.06'
$ws.Range('E17').Value = '  -5.94%  '
$ws.Range('D18').Value = "'331.45"
$ws.Range('E18').Value = '  -2.14%  '
$ws.Range('D19').Value = "'4.27"
$ws.Range('D20').Value = "'10.04"
$ws.Range('E20').Value = '  -3.25%  '
$ws.Range('D21').Value = "'6.11"
$ws.Range('E21').Value = '  -3.81%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').Value = "'63.96"
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('E25').Value = '  +0.29%  '
$ws.Range('E26').Value = '  -4.52%  '
$ws.Range('D27').Value = '2.656.63'
$ws.Range('E27').Value = '  -4.18%  '
$ws.Range('D28').Value = "'6.91"
$ws.Range('E28').Value = '  -2.38%  '
$ws.Range('D29').Value = '0.0₃0748'
$ws.Range('E29').Value = '  -6.60%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('E31').Value = '  -6.10%  '
$ws.Range('E32').Value = '  -2.95%  '
$ws.Range('D33').Value = "'148.47"
$ws.Range('E33').Value = '  -1.05%  '
$ws.Range('D34').Value = "'18.43"
$ws.Range('E34').Value = '  -2.24%  '
$ws.Range('D35').Value = "'3.95"
$ws.Range('E35').Value = '  -4.80%  '
$ws.Range('D36').Value = "'1.13"
$ws.Range('E36').Value = '  -5.43%  '
$ws.Range('D37').Value = "'0.840"
$ws.Range('E37').Value = '  -5.58%  '
$ws.Range('D38').Value = "'35.71"
$ws.Range('E38').Value = '  -2.99%  '
$ws.Range('D39').Value = "'0.817"
$ws.Range('E39').Value = '  -6.05%  '
$ws.Range('E40').Value = '  -3.95%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('D42').Value = "'3.46"
$ws.Range('E42').Value = '  -3.75%  '
$ws.Range('D43').Value = "'0.0951"
$ws.Range('E43').Value = '  -1.63%  '
$ws.Range('E44').Value = '  -0.71%  '
$ws.Range('D45').Value = "'0.573"
$ws.Range('E45').Value = '  -7.03%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = "'258.11"
$ws.Range('E46').Value = '  -6.33%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').Value = "'0.0519"
$ws.Range('E47').Value = '  -2.19%  '
$ws.Range('D48').Value = "'18.35"
$ws.Range('E48').Value = '  -7.70%  '
$ws.Range('D49').Value = '1.964.37'
$ws.Range('E49').Value = '  -4.09%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = "'0.0221"
$ws.Range('E50').Value = '  -3.33%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = "'4.50"
$ws.Range('E51').Value = '  -4.53%  '
